$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are decimal-looking text (e.g. "603.46") that Excel's COM
# layer would otherwise auto-coerce to a numeric type, losing the exact
# printed representation (trailing zeros, multi-dot grouping, etc). Forcing
# a text number format before the write keeps it a string; resetting the
# cell style back to "Normal" afterwards avoids leaving a visible style
# index on the cell (matches the original, unstyled cells).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '70.151.21'
$ws.Cells.Item(2, 5).Value = '  +0.35%  '
Set-TextValue $ws.Cells.Item(3, 4) '3.606.27'
$ws.Cells.Item(3, 5).Value = '  +2.87%  '
Set-TextValue $ws.Cells.Item(5, 4) '603.46'
$ws.Cells.Item(5, 5).Value = '  +0.63%  '
Set-TextValue $ws.Cells.Item(6, 4) '195.97'
$ws.Cells.Item(6, 5).Value = '  -0.03%  '
$ws.Cells.Item(7, 5).Value = '  +0.44%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.208'
$ws.Cells.Item(9, 5).Value = '  -0.51%  '
$ws.Cells.Item(10, 5).Value = '  -0.65%  '
Set-TextValue $ws.Cells.Item(11, 4) '53.81'
$ws.Cells.Item(11, 5).Value = '  -0.46%  '
$ws.Cells.Item(12, 5).Value = '  +1.52%  '
Set-TextValue $ws.Cells.Item(13, 4) '9.58'
$ws.Cells.Item(13, 5).Value = '  +0.42%  '
Set-TextValue $ws.Cells.Item(14, 4) '4.178.29'
$ws.Cells.Item(14, 5).Value = '  +2.93%  '
Set-TextValue $ws.Cells.Item(15, 4) '13.12'
$ws.Cells.Item(15, 5).Value = '  +4.19%  '
Set-TextValue $ws.Cells.Item(16, 4) '592.00'
$ws.Cells.Item(16, 5).Value = '  -2.39%  '
Set-TextValue $ws.Cells.Item(17, 4) '70.292.04'
$ws.Cells.Item(17, 5).Value = '  +0.36%  '
Set-TextValue $ws.Cells.Item(18, 4) '19.17'
$ws.Cells.Item(18, 5).Value = '  +0.99%  '
Set-TextValue $ws.Cells.Item(19, 4) '3.600.92'
$ws.Cells.Item(19, 5).Value = '  +2.57%  '
$ws.Cells.Item(20, 5).Value = '  +1.49%  '
Set-TextValue $ws.Cells.Item(21, 4) '0.995'
$ws.Cells.Item(21, 5).Value = '  +0.30%  '
Set-TextValue $ws.Cells.Item(22, 4) '17.69'
$ws.Cells.Item(22, 5).Value = '  -1.36%  '
Set-TextValue $ws.Cells.Item(23, 4) '5.17'
$ws.Cells.Item(23, 5).Value = '  +1.21%  '
Set-TextValue $ws.Cells.Item(24, 4) '101.96'
$ws.Cells.Item(24, 5).Value = '  -2.24%  '
$ws.Cells.Item(25, 5).Value = '  +0.45%  '
Set-TextValue $ws.Cells.Item(26, 4) '3.03'
$ws.Cells.Item(26, 5).Value = '  -1.01%  '
Set-TextValue $ws.Cells.Item(27, 4) '10.80'
$ws.Cells.Item(27, 5).Value = '  -1.66%  '
Set-TextValue $ws.Cells.Item(28, 4) '9.60'
$ws.Cells.Item(28, 5).Value = '  -0.96%  '
Set-TextValue $ws.Cells.Item(29, 4) '33.98'
$ws.Cells.Item(29, 5).Value = '  +1.24%  '
$ws.Cells.Item(30, 5).Value = '  +4.23%  '
Set-TextValue $ws.Cells.Item(31, 4) '7.12'
$ws.Cells.Item(31, 5).Value = '  +0.36%  '
Set-TextValue $ws.Cells.Item(32, 4) '12.31'
$ws.Cells.Item(32, 5).Value = '  -2.61%  '
$ws.Cells.Item(33, 5).Value = '  +1.15%  '
Set-TextValue $ws.Cells.Item(34, 4) '63.30'
$ws.Cells.Item(34, 5).Value = '  +0.21%  '
Set-TextValue $ws.Cells.Item(35, 4) '0.0₃0896'
$ws.Cells.Item(35, 5).Value = '  +8.99%  '
Set-TextValue $ws.Cells.Item(36, 4) '3.938.08'
$ws.Cells.Item(36, 5).Value = '  +5.67%  '
$ws.Cells.Item(37, 5).Value = '  +1.25%  '
Set-TextValue $ws.Cells.Item(38, 4) '528.08'
$ws.Cells.Item(38, 5).Value = '  +5.27%  '
$ws.Cells.Item(39, 5).Value = '  +0.06%  '
Set-TextValue $ws.Cells.Item(40, 4) '37.08'
$ws.Cells.Item(40, 5).Value = '  +0.74%  '
Set-TextValue $ws.Cells.Item(41, 4) '0.391'
$ws.Cells.Item(41, 5).Value = '  -0.47%  '
$ws.Cells.Item(42, 5).Value = '  -1.01%  '
$ws.Cells.Item(43, 5).Value = '  -1.42%  '
$ws.Cells.Item(44, 5).Value = '  -0.41%  '
Set-TextValue $ws.Cells.Item(45, 4) '3.43'
$ws.Cells.Item(45, 5).Value = '  +2.94%  '
Set-TextValue $ws.Cells.Item(46, 4) '2.86'
$ws.Cells.Item(46, 5).Value = '  +0.99%  '
$ws.Cells.Item(47, 5).Value = '  +0.64%  '
Set-TextValue $ws.Cells.Item(48, 4) '8.60'
$ws.Cells.Item(48, 5).Value = '  -1.33%  '
$ws.Cells.Item(49, 5).Value = '  -0.25%  '
$ws.Cells.Item(50, 5).Value = '  +4.94%  '
$ws.Cells.Item(51, 5).Value = '  +3.45%  '
